$p = $ppt.ActivePresentation

# --- Notes Master: datetimeFigureOut field (id {84B86612-B127-4CD0-BF15-66D49A7175F7}) ---
# 29-08-2025 -> 02-09-2025
$nm = $p.NotesMaster
$nm.Shapes.Item(2).TextFrame.TextRange.Text = "02-09-2025"

# --- Slide Master: datetimeFigureOut field (id {1D8BD707-D9CF-40AE-B4C6-C98DA3205C09}) ---
# 8/29/2025 -> 9/2/2025
$p.SlideMaster.Shapes.Item(14).TextFrame.TextRange.Text = "9/2/2025"

# --- Slide Layouts: same datetimeFigureOut field on each of the 5 layouts ---
$layouts = $p.SlideMaster.CustomLayouts
$layouts.Item(1).Shapes.Item(4).TextFrame.TextRange.Text = "9/2/2025"   # Title Slide
$layouts.Item(2).Shapes.Item(4).TextFrame.TextRange.Text = "9/2/2025"   # Title and Content
$layouts.Item(3).Shapes.Item(5).TextFrame.TextRange.Text = "9/2/2025"   # Two Content
$layouts.Item(4).Shapes.Item(3).TextFrame.TextRange.Text = "9/2/2025"   # Title Only
$layouts.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = "9/2/2025"   # Blank

# --- Slide 1: append NMID suffix to the "REGISTER NO AND NMID" line ---
$s1 = $p.Slides.Item(1)
$infoBox = $s1.Shapes.Item(7)
$infoBox.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "REGISTER NO AND NMID: 222407871/autunm1455ds242711"
